# MSME Country Indicators - Montenegro Summary
# Rename sheet "Data" -> "Summary", and add the Source Type / data row /
# source note below the existing Name + Title rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "Data" to "Summary"
$ws.Name = "Summary"

# 1a. Re-assert the pre-existing "name" (size 18) and "title" (bold)
#     formatting on the untouched header cells so it survives the
#     save round-trip unchanged.
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true

# 1b. Clear the old header row (5) and old data-label row (6); their
#     content is being replaced by the new rows 9/10 below.
$ws.Range("A5:D6").Clear()

# 2. New "Source Type" line (bold + underline) at row 7
$ws.Range("A7").Value = "Source Type: Statistical Institution"
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Underline = $true

# 3. Column headers (bold) move down to row 9
$ws.Range("B9").Value = "Micro"
$ws.Range("B9").Font.Bold = $true

$ws.Range("C9").Value = "SMEs"
$ws.Range("C9").Font.Bold = $true

$ws.Range("D9").Value = "MSMEs"
$ws.Range("D9").Font.Bold = $true

# 4. Data row at row 10: label stays bold, values are plain text
$ws.Range("A10").Value = "Enterprises (% of total)"
$ws.Range("A10").Font.Bold = $true

$ws.Range("B10").Value = "'30.8"
$ws.Range("C10").Value = "'64.9"
$ws.Range("D10").Value = "'95.7"

# 5. Source note (italic) at row 11
$ws.Range("A11").Value = "Source: SOM - Stat. Office of Montenegro, 2010"
$ws.Range("A11").Font.Italic = $true
